# "Generate Report for Handoff"
#
# The localization-status report records, per target language, the most
# recent handoff timestamp for each source file. A new handoff run was
# generated for the "4069360f-04e0-4a68-9b80-c4382b923214" source file
# (report row 5), so its "Latest Handoff Datetime" column (D) is updated
# on both the zh-cn and de-de status sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D5").Value = "2016-03-09 18:33:39"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D5").Value = "2016-03-09 18:33:44"
